$d = $word.ActiveDocument

# Find the 1-based index of the paragraph containing "LOB1039...".
$idx = 0
$targetIdx = -1
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -like "*LOB1039: Física Experimental III (Requisito fraco)*") {
        $targetIdx = $idx
        break
    }
}

if ($targetIdx -gt 0) {
    # Delete the three paragraphs following it: the blank line, the
    # "Ver no Jupiter..." line, and the "(c) 2020..." line.
    $pStart = $d.Paragraphs.Item($targetIdx + 1)
    $pEnd = $d.Paragraphs.Item($targetIdx + 3)

    $r = $d.Range($pStart.Range.Start, $pEnd.Range.End)
    $r.Delete()
}
